# Add a new textbox shape ("TextBox 4", id=5) to slide 1 with the repo link,
# matching the target OOXML:
#   <a:off x="7608699" y="108534"/><a:ext cx="4490937" cy="307777"/>
#   <a:bodyPr wrap="square"><a:spAutoFit/></a:bodyPr>
#   text run: lang="ko-KR" altLang="en-US" sz="1400" dirty="0"
#   "https://github.com/syesung01-max/IOT_week4_02.git"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoTextOrientationHorizontal = 1
$ppAutoSizeShapeToFitText = 1
$EMUsPerPoint = 12700

# The target shape lands as id=5 / "TextBox 4" -- i.e. the *second* textbox
# minted on this slide. Add + discard a throwaway one first so the
# id/name counter lands on the right value for the real shape.
$throwaway = $s.Shapes.AddTextbox($msoTextOrientationHorizontal, 0, 0, 10, 10)
$throwaway.Delete()

$left   = 7608699 / $EMUsPerPoint
$top    = 108534 / $EMUsPerPoint
$width  = 4490937 / $EMUsPerPoint
$height = 307777 / $EMUsPerPoint

$tb = $s.Shapes.AddTextbox($msoTextOrientationHorizontal, $left, $top, $width, $height)
$tb.Fill.Visible = $false

$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = $ppAutoSizeShapeToFitText

$tr = $tb.TextFrame.TextRange
$tr.Text = "https://github.com/syesung01-max/IOT_week4_02.git"
$tr.LanguageID = "ko-KR"
$tr.Font.Size = 14
